$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting A:E to B:F
$ws.Range("A:A").Insert()

# Fill in the sequence numbers for data rows 2-22 (21 rows) in the new column A
for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
